# Atualização de banco de dados
# Replace the user "evaldo" with "vanessa" throughout the "grants por usuario"
# sheet (3rd sheet), and tweak row 1 height / column B width / selection to
# match the author's final view state.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(3)

# Column B (rows 1-35) holds the username used by every GRANT/CREATE USER
# statement built in column D. Swap it from "evaldo" to "vanessa".
$ws.Range("B1:B35").Value = "vanessa"

# Row 1 grew taller and column B grew wider in the saved file (to fit the
# longer "vanessa" text), and the active selection moved to D1.
$ws.Range("1:1").RowHeight = 24
$ws.Columns.Item(2).ColumnWidth = 12.3

$ws.Range("D1").Select()
